# Refresh the cryptos Price/Volume(1h) columns (D:E) to match the
# latest scrape, per the GitHub Actions commit.
#
# D/E are plain text cells (quoted strings like "19.50" or
# "1.633.23", not real numbers), so any value that *looks* like a
# clean float is written with its NumberFormat forced to Text ("@")
# first -- otherwise Excel would coerce it to a number on assignment
# and silently drop significant trailing zeros (e.g. "19.50" -> 19.5).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.926.00'
$ws.Range("E2").Value = '  -0.27%  '
$ws.Range("D3").Value = '1.633.23'
$ws.Range("E3").Value = '  -0.43%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.94'
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5112'
$ws.Range("E6").Value = '  +0.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2577'
$ws.Range("E8").Value = '  +0.58%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06345'
$ws.Range("E9").Value = '  -0.34%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.50'
$ws.Range("E10").Value = '  -0.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07777'
$ws.Range("E11").Value = '  +0.21%  '
$ws.Range("E12").Value = '  -0.15%  '
$ws.Range("D13").Value = '1.636.71'
$ws.Range("E13").Value = '  -0.29%  '
$ws.Range("D14").Value = '1.858.00'
$ws.Range("E14").Value = '  -0.47%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5513'
$ws.Range("E15").Value = '  +1.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.95'
$ws.Range("E16").Value = '  -0.54%  '
$ws.Range("D17").Value = '0.0₅7648'
$ws.Range("D18").Value = '25.952.46'
$ws.Range("E18").Value = '  -0.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.003'
$ws.Range("E19").Value = '  +0.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '195.28'
$ws.Range("E20").Value = '  -0.75%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.420'
$ws.Range("E21").Value = '  -0.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.871'
$ws.Range("E22").Value = '  -0.58%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.051'
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.889'
$ws.Range("E25").Value = '  +0.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.92'
$ws.Range("E26").Value = '  +0.55%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1259'
$ws.Range("E27").Value = '  +5.51%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.60'
$ws.Range("E28").Value = '  -0.17%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.752'
$ws.Range("E29").Value = '  -1.40%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.241'
$ws.Range("E30").Value = '  +0.52%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04884'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.244'
$ws.Range("E32").Value = '  -0.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.192'
$ws.Range("E33").Value = '  +0.41%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.542'
$ws.Range("E34").Value = '  +0.74%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.372'
$ws.Range("E35").Value = '  +0.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.8973'
$ws.Range("E36").Value = '  +0.56%  '
$ws.Range("E37").Value = '  +1.70%  '
$ws.Range("E38").Value = '  -1.62%  '
$ws.Range("D39").Value = '1.117.08'
$ws.Range("E39").Value = '  -3.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01557'
$ws.Range("E40").Value = '  +0.10%  '
$ws.Range("E41").Value = '  +0.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.594'
$ws.Range("E42").Value = '  +2.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7960'
$ws.Range("E43").Value = '  -1.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.56'
$ws.Range("E44").Value = '  -1.40%  '
$ws.Range("D45").Value = '1.768.92'
$ws.Range("E45").Value = '  -0.44%  '
$ws.Range("E46").Value = '  -7.33%  '
$ws.Range("E47").Value = '  -1.65%  '
$ws.Range("E48").Value = '  +0.36%  '
$ws.Range("E49").Value = '  -0.18%  '
$ws.Range("E50").Value = '  +1.56%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.575'
$ws.Range("E51").Value = '  +3.27%  '
